$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1312.8422
$ws.Range("I100").Value = 1295.6923
$ws.Range("J100").Value = 1350
$ws.Range("K100").Value = 1295.6923
$ws.Range("L100").Value = 1350
$ws.Range("M100").Value = -754.6922999999999
$ws.Range("N100").Value = -2432

$ws.Range("H125").Value = 3751
$ws.Range("I125").Value = 866
$ws.Range("J125").Value = 4712.6665
$ws.Range("K125").Value = 7794
$ws.Range("L125").Value = 42413.9985
$ws.Range("M125").Value = -5334
$ws.Range("N125").Value = -47333.9985

$ws.Range("H138").Value = 3360.383
$ws.Range("I138").Value = 1494.909
$ws.Range("J138").Value = 4369.5737
$ws.Range("K138").Value = 4484.727000000001
$ws.Range("L138").Value = 13108.7211
$ws.Range("M138").Value = 655.2729999999992
$ws.Range("N138").Value = -23388.7211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 30900
$ws.Range("J123").Value = 30900
$ws.Range("L123").Value = 30900
$ws.Range("N123").Value = -40700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4160.7188
$ws.Range("I99").Value = 4812.4585
$ws.Range("J99").Value = 2205.5
$ws.Range("K99").Value = 4812.4585
$ws.Range("L99").Value = 2205.5
$ws.Range("M99").Value = -3314.4585
$ws.Range("N99").Value = -5201.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2143.275
$ws.Range("I31").Value = 1710.8077
$ws.Range("J31").Value = 2946.4285
$ws.Range("K31").Value = 1710.8077
$ws.Range("L31").Value = 2946.4285
$ws.Range("M31").Value = -1415.8077
$ws.Range("N31").Value = -3536.4285

$ws.Range("H34").Value = 2143.275
$ws.Range("I34").Value = 1710.8077
$ws.Range("J34").Value = 2946.4285
$ws.Range("K34").Value = 1710.8077
$ws.Range("L34").Value = 2946.4285
$ws.Range("M34").Value = -1508.8077
$ws.Range("N34").Value = -3350.4285

$ws.Range("H68").Value = 32283.334
$ws.Range("J68").Value = 32283.334
$ws.Range("L68").Value = 32283.334
$ws.Range("N68").Value = -33781.334

$ws.Range("H71").Value = 32283.334
$ws.Range("J71").Value = 32283.334
$ws.Range("L71").Value = 96850.00199999999
$ws.Range("N71").Value = -104338.002

$ws.Range("H107").Value = 914.8
$ws.Range("I107").Value = 771.53845
$ws.Range("J107").Value = 1180.8572
$ws.Range("K107").Value = 771.53845
$ws.Range("L107").Value = 1180.8572
$ws.Range("M107").Value = 1148.46155
$ws.Range("N107").Value = -5020.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2104.5789
$ws.Range("I69").Value = 1404.5
$ws.Range("J69").Value = 2427.6924
$ws.Range("K69").Value = 4213.5
$ws.Range("L69").Value = 7283.0772
$ws.Range("M69").Value = -3402.5
$ws.Range("N69").Value = -8905.0772

$ws.Range("H72").Value = 2104.5789
$ws.Range("I72").Value = 1404.5
$ws.Range("J72").Value = 2427.6924
$ws.Range("K72").Value = 12640.5
$ws.Range("L72").Value = 21849.2316
$ws.Range("M72").Value = -8584.5
$ws.Range("N72").Value = -29961.2316

$ws.Range("H80").Value = 1330
$ws.Range("I80").Value = 950
$ws.Range("J80").Value = 1482
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 4446
$ws.Range("M80").Value = -1914
$ws.Range("N80").Value = -6318

$ws.Range("H83").Value = 1330
$ws.Range("I83").Value = 950
$ws.Range("J83").Value = 1482
$ws.Range("K83").Value = 8550
$ws.Range("L83").Value = 13338
$ws.Range("M83").Value = -3870
$ws.Range("N83").Value = -22698

$ws.Range("H131").Value = 686.8387
$ws.Range("J131").Value = 829
$ws.Range("L131").Value = 2487
$ws.Range("N131").Value = -12567

$ws.Range("H132").Value = 1501
$ws.Range("I132").Value = 1002
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9018
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6488
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4763313.5
$ws.Range("I122").Value = 6251058.5
$ws.Range("K122").Value = 18753175.5
$ws.Range("M122").Value = -18750725.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16274

$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16288

$ws.Range("H59").Value = 9000
$ws.Range("J59").Value = 9000
$ws.Range("L59").Value = 9000
$ws.Range("N59").Value = -10308

$ws.Range("H68").Value = 3006.1538
$ws.Range("I68").Value = 3013.3333
$ws.Range("K68").Value = 3013.3333
$ws.Range("M68").Value = -2264.3333

$ws.Range("H71").Value = 3006.1538
$ws.Range("I71").Value = 3013.3333
$ws.Range("K71").Value = 15066.6665
$ws.Range("M71").Value = -11322.6665

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H87").Value = 32945
$ws.Range("J87").Value = 32945
$ws.Range("L87").Value = 32945
$ws.Range("N87").Value = -35191

$ws.Range("H88").Value = 24490
$ws.Range("J88").Value = 24490
$ws.Range("L88").Value = 24490
$ws.Range("N88").Value = -25346

$ws.Range("H90").Value = 32945
$ws.Range("J90").Value = 32945
$ws.Range("L90").Value = 98835
$ws.Range("N90").Value = -110067

$ws.Range("H91").Value = 24490
$ws.Range("J91").Value = 24490
$ws.Range("L91").Value = 24490
$ws.Range("N91").Value = -27454

$ws.Range("H93").Value = 11625.2
$ws.Range("I93").Value = 18875.666
$ws.Range("K93").Value = 18875.666
$ws.Range("M93").Value = -17627.666

$ws.Range("H111").Value = 35496.75
$ws.Range("J111").Value = 35496.75
$ws.Range("L111").Value = 35496.75
$ws.Range("N111").Value = -43676.75

$ws.Range("H112").Value = 37000
$ws.Range("J112").Value = 37000
$ws.Range("L112").Value = 37000
$ws.Range("N112").Value = -39954

$ws.Range("H123").Value = 38000
$ws.Range("J123").Value = 38000
$ws.Range("L123").Value = 38000
$ws.Range("N123").Value = -47800

$ws.Range("H128").Value = 55071.6
$ws.Range("J128").Value = 55071.6
$ws.Range("L128").Value = 55071.6
$ws.Range("N128").Value = -65031.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 48900
$ws.Range("J114").Value = 48900
$ws.Range("L114").Value = 48900
$ws.Range("N114").Value = -57578
